$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 5257.5
$ws.Range("I69").Value = 10000
$ws.Range("J69").Value = 3676.6667
$ws.Range("K69").Value = 30000
$ws.Range("L69").Value = 11030.0001
$ws.Range("M69").Value = -29126
$ws.Range("N69").Value = -12778.0001

$ws.Range("H72").Value = 5257.5
$ws.Range("I72").Value = 10000
$ws.Range("J72").Value = 3676.6667
$ws.Range("K72").Value = 90000
$ws.Range("L72").Value = 33090.0003
$ws.Range("M72").Value = -85632
$ws.Range("N72").Value = -41826.0003

$ws.Range("H113").Value = 3211.7693
$ws.Range("I113").Value = 1924.75
$ws.Range("K113").Value = 1924.75
$ws.Range("M113").Value = 1329.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1560.2
$ws.Range("I45").Value = 1646
$ws.Range("J45").Value = 1360
$ws.Range("K45").Value = 1646
$ws.Range("L45").Value = 1360
$ws.Range("M45").Value = -1269
$ws.Range("N45").Value = -2114

$ws.Range("H122").Value = 1664.7693
$ws.Range("I122").Value = 1728.5
$ws.Range("J122").Value = 1562.8
$ws.Range("K122").Value = 5185.5
$ws.Range("L122").Value = 4688.4
$ws.Range("M122").Value = -2735.5
$ws.Range("N122").Value = -9588.4

$ws.Range("H132").Value = 2744.842
$ws.Range("I132").Value = 2701.5
$ws.Range("J132").Value = 2976
$ws.Range("K132").Value = 8104.5
$ws.Range("L132").Value = 8928
$ws.Range("M132").Value = -5574.5
$ws.Range("N132").Value = -13988

$ws.Range("H139").Value = 152782
$ws.Range("J139").Value = 152782
$ws.Range("L139").Value = 152782
$ws.Range("N139").Value = -163062

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2255.9375
$ws.Range("I20").Value = 2551.25
$ws.Range("J20").Value = 1370
$ws.Range("K20").Value = 2551.25
$ws.Range("L20").Value = 1370
$ws.Range("M20").Value = -2304.25
$ws.Range("N20").Value = -1864

$ws.Range("H99").Value = 1512.5
$ws.Range("I99").Value = 1350
$ws.Range("K99").Value = 1350
$ws.Range("M99").Value = 148

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2848.3333
$ws.Range("I16").Value = 1822.5
$ws.Range("J16").Value = 4900
$ws.Range("K16").Value = 1822.5
$ws.Range("L16").Value = 4900
$ws.Range("M16").Value = -1535.5
$ws.Range("N16").Value = -5474

$ws.Range("H99").Value = 3429.1428
$ws.Range("I99").Value = 2376
$ws.Range("J99").Value = 4833.3335
$ws.Range("K99").Value = 2376
$ws.Range("L99").Value = 4833.3335
$ws.Range("M99").Value = -878
$ws.Range("N99").Value = -7829.3335

$ws.Range("H113").Value = 2848.3333
$ws.Range("I113").Value = 1822.5
$ws.Range("J113").Value = 4900
$ws.Range("K113").Value = 1822.5
$ws.Range("L113").Value = 4900
$ws.Range("M113").Value = 347.5
$ws.Range("N113").Value = -9240

$ws.Range("H126").Value = 3429.1428
$ws.Range("I126").Value = 2376
$ws.Range("J126").Value = 4833.3335
$ws.Range("K126").Value = 7128
$ws.Range("L126").Value = 14500.0005
$ws.Range("M126").Value = -4658
$ws.Range("N126").Value = -19440.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3047.8823
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 3781.4
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 11344.2
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -16284.2

$ws.Range("H132").Value = 22750.092
$ws.Range("I132").Value = 29756
$ws.Range("K132").Value = 89268
$ws.Range("M132").Value = -86738

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3020.75
$ws.Range("I7").Value = 3178.8
$ws.Range("J7").Value = 650
$ws.Range("K7").Value = 3178.8
$ws.Range("L7").Value = 650
$ws.Range("M7").Value = -3066.8
$ws.Range("N7").Value = -874

$ws.Range("H40").Value = 3121.6667
$ws.Range("I40").Value = 2994.4443
$ws.Range("J40").Value = 3503.3333
$ws.Range("K40").Value = 2994.4443
$ws.Range("L40").Value = 3503.3333
$ws.Range("M40").Value = -2858.4443
$ws.Range("N40").Value = -3775.3333

$ws.Range("H46").Value = 1142.8572
$ws.Range("I46").Value = 1466.6666
$ws.Range("J46").Value = 900
$ws.Range("K46").Value = 1466.6666
$ws.Range("L46").Value = 900
$ws.Range("M46").Value = -1278.6666
$ws.Range("N46").Value = -1276

$ws.Range("H68").Value = 3750
$ws.Range("I68").Value = 3333.3333
$ws.Range("J68").Value = 5000
$ws.Range("K68").Value = 3333.3333
$ws.Range("L68").Value = 5000
$ws.Range("M68").Value = -2584.3333
$ws.Range("N68").Value = -6498

$ws.Range("H71").Value = 3750
$ws.Range("I71").Value = 3333.3333
$ws.Range("J71").Value = 5000
$ws.Range("K71").Value = 16666.6665
$ws.Range("L71").Value = 25000
$ws.Range("M71").Value = -12922.6665
$ws.Range("N71").Value = -32488

$ws.Range("H122").Value = 8072.1787
$ws.Range("I122").Value = 7736.2354
$ws.Range("J122").Value = 8591.362999999999
$ws.Range("K122").Value = 23208.7062
$ws.Range("L122").Value = 25774.089
$ws.Range("M122").Value = -20758.7062
$ws.Range("N122").Value = -30674.089

$ws.Range("H126").Value = 3020.75
$ws.Range("I126").Value = 3178.8
$ws.Range("J126").Value = 650
$ws.Range("K126").Value = 9536.400000000001
$ws.Range("L126").Value = 1950
$ws.Range("M126").Value = -7066.400000000001
$ws.Range("N126").Value = -6890

$ws.Range("H132").Value = 6200.7
$ws.Range("I132").Value = 7668
$ws.Range("J132").Value = 5571.857
$ws.Range("K132").Value = 23004
$ws.Range("L132").Value = 16715.571
$ws.Range("M132").Value = -20474
$ws.Range("N132").Value = -21775.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3546.5454
$ws.Range("I62").Value = 3583.3333
$ws.Range("J62").Value = 3502.4
$ws.Range("K62").Value = 3583.3333
$ws.Range("L62").Value = 3502.4
$ws.Range("M62").Value = -2959.3333
$ws.Range("N62").Value = -4750.4

$ws.Range("H65").Value = 3546.5454
$ws.Range("I65").Value = 3583.3333
$ws.Range("J65").Value = 3502.4
$ws.Range("K65").Value = 17916.6665
$ws.Range("L65").Value = 17512
$ws.Range("M65").Value = -14796.6665
$ws.Range("N65").Value = -23752

$ws.Range("H122").Value = 5161.6787
$ws.Range("I122").Value = 1566.5
$ws.Range("J122").Value = 14149.625
$ws.Range("K122").Value = 4699.5
$ws.Range("L122").Value = 42448.875
$ws.Range("M122").Value = -2249.5
$ws.Range("N122").Value = -47348.875

$ws.Range("H126").Value = 1789.0588
$ws.Range("I126").Value = 1942
$ws.Range("J126").Value = 1422
$ws.Range("K126").Value = 5826
$ws.Range("L126").Value = 4266
$ws.Range("M126").Value = -3356
$ws.Range("N126").Value = -9206

$ws.Range("H132").Value = 4559.2354
$ws.Range("I132").Value = 3846.4614
$ws.Range("K132").Value = 11539.3842
$ws.Range("M132").Value = -9009.3842
